$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 7 (the "ctrl 1.2.3.txt" row), shifting rows below it up
$ws.Rows("7:7").Delete()

# Set a custom width for column A (closest achievable value to the
# target stored width of 20.3671875 given this runtime's width quantization)
$ws.Columns("A").ColumnWidth = 19.5

# Change the selection to a single cell D5 (instead of the whole A1:B85 range)
$ws.Range("D5").Select()
